# Fruta / hortaliza, semanal
# Insert 4 new weekly data rows into the Nectarín sheet (before what was row 876),
# shifting all subsequent rows down by 4 (954 -> 958 total data rows),
# and populate the newly inserted rows with the new weekly price-report entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 876..954 down to 880..958
$ws.Rows("876:879").Insert()

$newRows = @(
    @(10, "Vega Modelo de Temuco", "La Araucanía", 45265, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Early Glo",  "Primera", 2,   580000, 580000, 580000, "$/bins (420 kilos)",        "Región de O'Higgins", 1381, 420),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 45265, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Early John", "Primera", 250, 25000,  25000,  25000,  "$/bandeja 18 kilos granel", "Región de O'Higgins", 1389, 18),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 45265, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Super Queen","Primera", 300, 22000,  22000,  22000,  "$/bandeja 18 kilos granel", "Región de O'Higgins", 1222, 18),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 45265, 9, "Fruta", 100103, "Frutos de hueso (carozo)", 100103006, "Nectarín", "Super Queen","Primera", 3,   480000, 480000, 480000, "$/bins (420 kilos)",        "Región de O'Higgins", 1143, 420)
)

$startRow = 876
for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = $startRow + $i
    $rowValues = $newRows[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowValues[$c]
    }
}
